$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Duplicate value in FA18 column" testdata file references to FA19.
$ws.Range("C16").Value = "ICER RRMM - Duplicate value in FA19 column.xlsx"
$ws.Range("D16").Value = "\Testdata\Templates\ImportPublications\Testing_Env\ICER RRMM - Duplicate value in FA19 column.xlsx"

# Update "ICER - ICER RRMM 2022 report - 12/19/2022" -> "ICER - ICER RRMM 2022 report"
# (drop the trailing date) in every cell of column B that references it.
$newReportName = "ICER - ICER RRMM 2022 report"
$ws.Range("B2").Value = $newReportName
$ws.Range("B4").Value = $newReportName
$ws.Range("B8").Value = $newReportName
$ws.Range("B12").Value = $newReportName
$ws.Range("B16").Value = $newReportName

# Adjust the stored view: no frozen/topmost left column, selection moved to C9.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C9").Select() | Out-Null

# Column B width changed (narrower) to fit the shorter report name text.
$ws.Range("B:B").ColumnWidth = 25.67
